$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume (E) columns keep their text representation
# (values like "1.000", "0.00001087" etc. must not be auto-converted to numbers)
$ws.Range("D2:E51").NumberFormat = "@"

# Rows 2-43: simple D/E value updates
$ws.Range("D2").Value = "29.454.31"
$ws.Range("D3").Value = "1.853.12"
$ws.Range("E3").Value = "  +0.38%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "240.84"
$ws.Range("E5").Value = "  +0.88%  "
$ws.Range("D6").Value = "0.6311"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "0.07675"
$ws.Range("E8").Value = "  +1.68%  "
$ws.Range("D9").Value = "0.2945"
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("E10").Value = "  +0.22%  "
$ws.Range("D11").Value = "0.07754"
$ws.Range("E11").Value = "  +0.69%  "
$ws.Range("D12").Value = "1.851.62"
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("E13").Value = "  +0.93%  "
$ws.Range("D14").Value = "0.00001087"
$ws.Range("E14").Value = "  +6.68%  "
$ws.Range("D15").Value = "0.6819"
$ws.Range("E15").Value = "  +0.25%  "
$ws.Range("D16").Value = "83.65"
$ws.Range("E16").Value = "  +0.61%  "
$ws.Range("D17").Value = "2.109.05"
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("D18").Value = "6.171"
$ws.Range("E18").Value = "  +0.63%  "
$ws.Range("D19").Value = "29.473.85"
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("D20").Value = "229.63"
$ws.Range("E20").Value = "  +0.60%  "
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").Value = "7.458"
$ws.Range("E23").Value = "  -0.46%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").Value = "156.90"
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("D26").Value = "0.1386"
$ws.Range("E26").Value = "  -0.74%  "
$ws.Range("D27").Value = "8.413"
$ws.Range("E27").Value = "  +0.49%  "
$ws.Range("D28").Value = "17.73"
$ws.Range("E28").Value = "  +0.38%  "
$ws.Range("D29").Value = "1.327"
$ws.Range("E29").Value = "  +4.34%  "
$ws.Range("D30").Value = "1.473"
$ws.Range("E30").Value = "  +0.92%  "
$ws.Range("D31").Value = "0.05690"
$ws.Range("E31").Value = "  +0.56%  "
$ws.Range("D32").Value = "4.133"
$ws.Range("E32").Value = "  +0.23%  "
$ws.Range("D33").Value = "4.053"
$ws.Range("E33").Value = "  +0.54%  "
$ws.Range("D34").Value = "1.852"
$ws.Range("E34").Value = "  +0.68%  "
$ws.Range("D35").Value = "1.165"
$ws.Range("E35").Value = "  +0.66%  "
$ws.Range("D36").Value = "0.7053"
$ws.Range("E36").Value = "  -1.43%  "
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("D38").Value = "2.784"
$ws.Range("E38").Value = "  +0.36%  "
$ws.Range("E39").Value = "  -0.57%  "
$ws.Range("D40").Value = "1.218.72"
$ws.Range("E40").Value = "  -2.24%  "
$ws.Range("D41").Value = "6.553"
$ws.Range("E41").Value = "  +5.87%  "
$ws.Range("D42").Value = "0.9103"
$ws.Range("E42").Value = "  +0.84%  "
$ws.Range("E43").Value = "  +0.07%  "

# Rows 44-51: shifted block (new RocketPoolETH row inserted at 44, Algorand dropped)
$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").Value = "2.017.63"
$ws.Range("E44").Value = "  +0.16%  "
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "101.72"
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "66.50"
$ws.Range("E46").Value = "  +0.21%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.00000000121"
$ws.Range("E47").Value = "  +0.98%  "
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").Value = "7.126"
$ws.Range("E48").Value = "  +0.68%  "
$ws.Range("B49").Value = "TheSandbox"
$ws.Range("C49").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D49").Value = "0.4021"
$ws.Range("E49").Value = "  +0.22%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "9.043"
$ws.Range("E50").Value = "  +0.19%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").Value = "1.686"
$ws.Range("E51").Value = "  -0.32%  "
